$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.287.57"
$ws.Range("E2").Value = "  -4.80%  "
$ws.Range("D3").Value = "2.906.19"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.85"
$ws.Range("E5").Value = "  -3.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.78"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").Value = "2.900.91"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -6.20%  "
$ws.Range("E11").Value = "  -6.60%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000213"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.11"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.383.94"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "2.899.53"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.55"
$ws.Range("E18").Value = "  +7.00%  "
$ws.Range("D19").Value = "57.249.38"
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "405.51"
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.83"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.01"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.68"
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0977"
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.916"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.02"
$ws.Range("E36").Value = "  -9.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.98"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.21"
$ws.Range("E38").Value = "  +5.14%  "
$ws.Range("D39").Value = "0.0₃0633"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0339"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.45"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "364.86"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "2.610.44"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.88"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.108"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.229"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.58"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -1.08%  "
